# "VENTAS POR GRUPO" gains a new product-category column "GRANITO" (inserted
# right before "GRIFERIAS", i.e. at column F, pushing the existing F:N
# columns one place right to G:O), plus three brand-new trailing columns:
# "NO RESURTIBLES", "PANELES PVC" and "PANELES PU" (P:R).
#
# Note: Excel's ColumnWidth COM property is in "characters" and is offset
# from the width actually persisted in the OOXML <col width="..."> attribute
# by a constant 5/6 (≈0.8333) of a character (the default-font grid-line
# padding). Subtracting 5/6 from the target stored width reproduces the
# exact width values from the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert "GRANITO" before "GRIFERIAS" (new column F) ---------------
$ws.Columns("F:F").Insert()

$ws.Range("F1").Value = "GRANITO"
$ws.Range("F2").Value = 0
$ws.Range("F3").Value = "0 de 1"

$ws.Columns("F:F").ColumnWidth = (13 - 5/6)

# --- Append three new trailing columns: P, Q, R ------------------------
# Clone the formatting (header/number/"x de y" styles) from the last
# existing column (O, "SAL SOLUBLE") onto the new P:R columns first …
$ws.Range("O1:O3").Copy()
$ws.Range("P1:R3").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# … then fill in the actual header/row values.
$ws.Range("P1").Value = "NO RESURTIBLES"
$ws.Range("Q1").Value = "PANELES PVC"
$ws.Range("R1").Value = "PANELES PU"

$ws.Range("P2:R2").Value = 0
$ws.Range("P3:R3").Value = "0 de 1"

$ws.Columns("P:P").ColumnWidth = (20 - 5/6)
$ws.Columns("Q:Q").ColumnWidth = (17 - 5/6)
$ws.Columns("R:R").ColumnWidth = (16 - 5/6)
